$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54, shifting existing rows 54-106 down to 55-107.
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new weekly data point.
$ws.Range("A54").Value = 8
$ws.Range("B54").Value = "Terminal La Palmera de La Serena"
$ws.Range("C54").Value = "Coquimbo"
$ws.Range("D54").Value = 44554
$ws.Range("E54").Value = 4
$ws.Range("F54").Value = 100112001
$ws.Range("G54").Value = "Berenjena"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 500
$ws.Range("K54").Value = 10000
$ws.Range("L54").Value = 11000
$ws.Range("M54").Value = 10500
$ws.Range("N54").Value = "$/caja 60 unidades"
$ws.Range("O54").Value = "Región de Arica y Parinacota"
$ws.Range("P54").Value = 175
$ws.Range("Q54").Value = 60
$ws.Range("R54").Value = "Hortaliza"
